$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Docentes responsáveis:" / "Semestral") is removed entirely; this
# shifts rows 13-26 up to become rows 12-25.
$ws.Rows.Item(12).Delete()

# Row 10 ("Objetivos:") gets a new long-form description in B/C.
$ws.Range("B10").Value = "A ciência dos biomateriais é uma atividade multidisciplinar que envolve a medicina, as ciências naturais e as engenharias, delimitando duas grandes áreas: a biotecnologia e a bioengenharia. A disciplina Biomateriais visa prover aos estudantes fundamentos básicos da ciência de biomateriais, dar uma perspectiva sobre os principais biomateriais aplicados em algumas áreas da medicina e contribuir para a compreensão das interações célula-material. Dessa forma, contribuir para o desenvolvimento da área e certamente alavancar a formação de recursos humanos associados a um melhor uso da infra-estrutura já existente."
$ws.Range("C10").Value = "A ciência dos biomateriais é uma atividade multidisciplinar que envolve a medicina, as ciências naturais e as engenharias, delimitando duas grandes áreas: a biotecnologia e a bioengenharia. A disciplina Biomateriais visa prover aos estudantes fundamentos básicos da ciência de biomateriais, dar uma perspectiva sobre os principais biomateriais aplicados em algumas áreas da medicina e contribuir para a compreensão das interações célula-material. Dessa forma, contribuir para o desenvolvimento da área e certamente alavancar a formação de recursos humanos associados a um melhor uso da infra-estrutura já existente."

# Row 12 (post-shift, was row 13: "Programa resumido:") gets the
# "2166002 - Sandra Giacomin Schneider" value that used to sit under
# "Objetivos:" before the row shift.
$ws.Range("B12").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C12").Value = "2166002 - Sandra Giacomin Schneider"

# Row 14 (post-shift, was row 15: "Programa:") gets the short syllabus text.
$ws.Range("B14").Value = "1 - Introdução aos Biomateriais`n2 - Interação tecido - implante`n3 - Técnicas de modificação de superfície`n4 - Técnicas de caracterização biológica`n5 - Aspectos práticos no uso de biomateriais"
$ws.Range("C14").Value = "1 - Introdução aos Biomateriais`n2 - Interação tecido - implante`n3 - Técnicas de modificação de superfície`n4 - Técnicas de caracterização biológica`n5 - Aspectos práticos no uso de biomateriais"

# Row 17 (post-shift, was row 18: "Método:") gets the detailed method text.
$ws.Range("B17").Value = "1 - Introdução aos Biomateriais`n  1.1- Conceitos básicos de biomateriais; `n  1.2 - Classes de materiais usados na área biomédica;`n  1.3 - Classificação dos biomateriais quanto à resposta biológica`n2 - Interação tecido  implante:`n  2.1 - Histórico da osteointegração; `n  2.2 - Fisiologia do osso;`n  2.3 - Natureza da ligação osso-implante;`n  2.4 - Aspectos superficiais dos implantes.`n3 - Técnicas de modificação da superfície:`n  3.1 - Técnicas para criar uma superfície bioativa: cerâmicas bioativas e biovidros, recobrimentos com fosfatos de  cálcio como transportador de proteínas ósseas morfogenéticas;`n  3.2 - Técnicas para aumentar a rugosidade superficial: usinagem, ataque ácido, jateamento, aspersão térmica. `n4 - Técnicas de caracterização biológica`n  4.1 - Teste em líquido corporal simulado (SBF)`n  4.2 - Cultura de células (in vitro)`n  4.3  Teste com cobaias (in vivo)`n5 - Aspectos práticos no uso de biomateriais`n  5.1- Técnicas de esterilização`n  5.2  Normas técnicas"
$ws.Range("C17").Value = "1 - Introdução aos Biomateriais`n  1.1- Conceitos básicos de biomateriais; `n  1.2 - Classes de materiais usados na área biomédica;`n  1.3 - Classificação dos biomateriais quanto à resposta biológica`n2 - Interação tecido  implante:`n  2.1 - Histórico da osteointegração; `n  2.2 - Fisiologia do osso;`n  2.3 - Natureza da ligação osso-implante;`n  2.4 - Aspectos superficiais dos implantes.`n3 - Técnicas de modificação da superfície:`n  3.1 - Técnicas para criar uma superfície bioativa: cerâmicas bioativas e biovidros, recobrimentos com fosfatos de  cálcio como transportador de proteínas ósseas morfogenéticas;`n  3.2 - Técnicas para aumentar a rugosidade superficial: usinagem, ataque ácido, jateamento, aspersão térmica. `n4 - Técnicas de caracterização biológica`n  4.1 - Teste em líquido corporal simulado (SBF)`n  4.2 - Cultura de células (in vitro)`n  4.3  Teste com cobaias (in vivo)`n5 - Aspectos práticos no uso de biomateriais`n  5.1- Técnicas de esterilização`n  5.2  Normas técnicas"
